# "utilizando densidade populacional para classificar grandes cidades"
#
# Reclassify the "ideologia" column (B):
#   - "Centro Dem" -> "Centro"
#   - "Centrao"    -> "Centro"
#   - "Extinto"    -> "Direita"   (PHS row)
# and correct a handful of the yearly seat-count figures that shift as a
# result of the new classification.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reclassify ideologia (column B) -------------------------------------
for ($r = 2; $r -le 26; $r++) {
    $ideologia = $ws.Cells.Item($r, 2).Value2
    if ($ideologia -eq "Centro Dem") {
        $ws.Cells.Item($r, 2).Value = "Centro"
    } elseif ($ideologia -eq "Centrao") {
        $ws.Cells.Item($r, 2).Value = "Centro"
    } elseif ($ideologia -eq "Extinto") {
        $ws.Cells.Item($r, 2).Value = "Direita"
    }
}

# --- Corrected seat counts -------------------------------------------------
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 6
$ws.Range("F4").Value = 24
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = 8
$ws.Range("F17").Value = 9
